$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("A8").NumberFormat = "@"
$ws.Range("A8").Value = ''
$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 40
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.0'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = 'Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it''s ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet''s & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F8").Value = 0
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '0.00'
$ws.Range("H8").Value = 0

# Row 9
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = 'P. point'
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 63
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = 'Short point (up to 3 mtr.)'
$ws.Range("F9").Value = 256
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '16128.00'
$ws.Range("H9").Value = 0

# Row 10
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = 'P. point'
$ws.Range("B10").Value = 0
$ws.Range("C10").Value = 1
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = 'Medium point (up to 6 mtr.)'
$ws.Range("F10").Value = 472
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '472.00'
$ws.Range("H10").Value = 0

# Row 11
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = 'P. point'
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 13
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '4'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = 'Long point  (up to 10 mtr.)'
$ws.Range("F11").Value = 662
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '8606.00'
$ws.Range("H11").Value = 0

# Row 12
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = 'Each'
$ws.Range("B12").Value = 0
$ws.Range("C12").Value = 82
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '3.0'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = 'P & F ISI marked (IS:3854) 6 amp. flush type non modular switch  with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F12").Value = 23
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '1886.00'
$ws.Range("H12").Value = 0

# Row 13
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = 'Each'
$ws.Range("B13").Value = 0
$ws.Range("C13").Value = 29
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.0'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = 'P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F13").Value = 50
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '1450.00'
$ws.Range("H13").Value = 0

# Row 14
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = 'Each'
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 93
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.0'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F14").Value = 78
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '7254.00'
$ws.Range("H14").Value = 0

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = 'Each'
$ws.Range("B15").Value = 0
$ws.Range("C15").Value = 76
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.0'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = 'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F15").Value = 30
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '2280.00'
$ws.Range("H15").Value = 0

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = 'Each'
$ws.Range("B16").Value = 0
$ws.Range("C16").Value = 15
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '9.0'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F16").Value = 219
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '3285.00'
$ws.Range("H16").Value = 0

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = 'Each'
$ws.Range("B17").Value = 0
$ws.Range("C17").Value = 24
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '10.0'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F17").Value = 303
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '7272.00'
$ws.Range("H17").Value = 0

# Row 18
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = 'R. mtr.'
$ws.Range("B18").Value = 72
$ws.Range("C18").Value = 72
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '17'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '25 mm'
$ws.Range("F18").Value = 56
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '4032.00'
$ws.Range("H18").Value = 0

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = 'Mtr.'
$ws.Range("B19").Value = 0
$ws.Range("C19").Value = 61
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '2 x 2.5 sq. mm. + 1x1.5sqmm'
$ws.Range("F19").Value = 81
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '4941.00'
$ws.Range("H19").Value = 0

# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = ''
$ws.Range("B20").Value = 0
$ws.Range("C20").Value = 93
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.0'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
$ws.Range("F20").Value = 0
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '0.00'
$ws.Range("H20").Value = 0

# Row 21
$ws.Range("A21").NumberFormat = "@"
$ws.Range("A21").Value = ''
$ws.Range("B21").Value = 0
$ws.Range("C21").Value = 76
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.0'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F21").Value = 0
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '0.00'
$ws.Range("H21").Value = 0

# Row 22
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = ''
$ws.Range("B22").Value = 0
$ws.Range("C22").Value = 44
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '29'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = 'Single pole MCB   (With B/C curve tripping Characteristics)'
$ws.Range("F22").Value = 0
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '0.00'
$ws.Range("H22").Value = 0

# Row 23
$ws.Range("B23").Value = 0
$ws.Range("C23").Value = 77
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '31'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = 'Double pole MCB(With B/C curve tripping Characteristics)'
$ws.Range("F23").Value = 0
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '0.00'
$ws.Range("H23").Value = 0

# Row 24
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = 'Each'
$ws.Range("B24").Value = 0
$ws.Range("C24").Value = 77
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '32'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = ' 50/63 A rating'
$ws.Range("F24").Value = 900
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '69300.00'
$ws.Range("H24").Value = 0

# Row 25
$ws.Range("B25").Value = 0
$ws.Range("C25").Value = 60
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '18.0'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = 'Providing & Fixing of Recessed/surface mounting heavy duty horizontal type Double Door ( Metal / Glazed )Distribution board with Metal end box made out from Galvanized steel / CRCA sheet not less then 1.2 mm thick  conforming to IS-8623-1 & 3 /  IEC 61439- 1 & 3, powder painted complete with reversible door (for double door DB only )100 amp.  insulated copper bus bar/shorting link , copper neutral link, copper earth link , color coded interconnecting wire set  of suitable rating and din bar,masking sheet,  making internal DB  terminations with copper lugs, Ferrules,  detachable gland plate, including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
$ws.Range("F25").Value = 0
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '0.00'
$ws.Range("H25").Value = 0

# Row 26
$ws.Range("A26").NumberFormat = "@"
$ws.Range("A26").Value = ''
$ws.Range("B26").Value = 0
$ws.Range("C26").Value = 80
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '34'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = 'Metal door (single phase) IK-09 and IP-43 with Metal end box'
$ws.Range("F26").Value = 0
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '0.00'
$ws.Range("H26").Value = 0

# Row 27
$ws.Range("B27").Value = 0
$ws.Range("C27").Value = 34
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '36'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = 'Total'
$ws.Range("F27").Value = 0
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '0.00'
$ws.Range("H27").Value = 0

# Row 28 (new content)
$ws.Range("A28").NumberFormat = "@"
$ws.Range("A28").Value = '%'
$ws.Range("B28").Value = 0
$ws.Range("C28").Value = 51
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '37'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = 'Add Tender Premium '
$ws.Range("F28").Value = 0
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '0.00'
$ws.Range("H28").Value = 0
$ws.Range("I28").NumberFormat = "@"
$ws.Range("I28").Value = ''

# Row 29 (clear Grand Total content, leave only blank row)
$ws.Range("E29").ClearContents()
$ws.Range("G29").ClearContents()
$ws.Range("H29").ClearContents()

# Row 30 (was Tender Premium, becomes Grand Total)
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = 'Grand Total Rs.'
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '126906.00'
$ws.Range("H30").NumberFormat = "@"
$ws.Range("H30").Value = '126906.00'

# Row 31 (was NET PAYABLE, becomes Tender Premium @ 0%)
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = 'Tender Premium @ 0%'
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '0.00'
$ws.Range("H31").NumberFormat = "@"
$ws.Range("H31").Value = '0.00'

# Row 32 (new NET PAYABLE AMOUNT row)
$ws.Range("A32").NumberFormat = "@"
$ws.Range("A32").Value = ''
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = ''
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = ''
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = ''
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = 'NET PAYABLE AMOUNT Rs.'
$ws.Range("F32").NumberFormat = "@"
$ws.Range("F32").Value = ''
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '126906.00'
$ws.Range("H32").NumberFormat = "@"
$ws.Range("H32").Value = '126906.00'
$ws.Range("I32").NumberFormat = "@"
$ws.Range("I32").Value = ''

Write-Output "Edit complete"
